$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Block 1: rows 2-21 (UnitMass column C), "+ loading" table
$ws.Range("C2").Value = 23
$ws.Range("C3").Value = 28
$ws.Range("C4").Value = 125
$ws.Range("C5").Value = 27
$ws.Range("C6").Value = 83
$ws.Range("C7").Value = 69
$ws.Range("C8").Value = 82
$ws.Range("C9").Value = 68
$ws.Range("C10").Value = 45
$ws.Range("C11").Value = 58
$ws.Range("C12").Value = 81
$ws.Range("C13").Value = 19
$ws.Range("C14").Value = 1
$ws.Range("C15").Value = 61
$ws.Range("C16").Value = 70
$ws.Range("C17").Value = 15
$ws.Range("C18").Value = 95
$ws.Range("C19").Value = 98
$ws.Range("C20").Value = 126
$ws.Range("C21").Value = 14

# Block 2: rows 23-42 (UnitMass column C), "- loading" table
$ws.Range("C23").Value = 39
$ws.Range("C24").Value = 71
$ws.Range("C25").Value = 175
$ws.Range("C26").Value = 55
$ws.Range("C27").Value = 41
$ws.Range("C28").Value = 43
$ws.Range("C29").Value = 231
$ws.Range("C30").Value = 115
$ws.Range("C31").Value = 102
$ws.Range("C32").Value = 103
$ws.Range("C33").Value = 149
$ws.Range("C34").Value = 91
$ws.Range("C35").Value = 59
$ws.Range("C36").Value = 112
$ws.Range("C37").Value = 46
$ws.Range("C38").Value = 7
$ws.Range("C39").Value = 287
$ws.Range("C40").Value = 77
$ws.Range("C41").Value = 88
$ws.Range("C42").Value = 159

$wb.Save()
